# Buttons workbook customization edit
# - Rename Sheet1 -> Buttons
# - Insert a new "QtDesigner" row (App) as row 4, pushing the rest of the
#   table down by one row
# - Remove the yellow highlight from the two "Downloads / Folder" rows
#   (now rows 7 and 12 after the insert)
# - Widen column C for the longer paths
# - Update the active selection to the new last data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Buttons"

# Insert a new row at position 4 (shifts old rows 4-12 down to 5-13)
$ws.Rows.Item(4).Insert() | Out-Null

$ws.Range("A4").Value = "QtDesigner"
$ws.Range("B4").Value = "App"
$ws.Range("C4").Value = "C:\Users\tuan\Anaconda3\envs\icrm\Scripts\designer.exe"
$ws.Range("D4").Value = "tuan"

# Drop the yellow "Downloads" highlight (it shifted from rows 6/11 to 7/12)
$ws.Range("A7:C7").Style = "Normal"
$ws.Range("A12:C12").Style = "Normal"

# Column C needs to be wider to fit the new / existing long paths
$ws.Columns.Item(3).ColumnWidth = 68.42

# Update view state to match: selection on what is now the next empty row
$ws.Range("C14").Select() | Out-Null

Write-Host "Edit applied"
